# New STH Limits module integrated everywhere STU Limits module new and integrated
#
# Updates a handful of data cells (column E) across the "System
# Configuration", "Statistics" and "Calibration" sheets. Several of the
# new values look numeric (dates / counters / floats) but must be written
# back as plain text, exactly like the existing cells in that column, so
# we route those through a scratch cell + Copy/PasteSpecial(Values) trick
# to avoid Excel's automatic text->number coercion on a direct .Value
# assignment (which would also tack on a spurious quote-prefix style).

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$cellAddr, [string]$text)

    $scratch = $ws.Cells.Item(500, 500)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.Clear()
    $excel.CutCopyMode = $false
}

# --- System Configuration0@0x0 ---
$wsSystem = $wb.Worksheets.Item("System Configuration0@0x0")
$wsSystem.Range("E3").Value = "D20-4889"

# --- Statistics@0x5 ---
$wsStats = $wb.Worksheets.Item("Statistics@0x5")
Set-TextValue $wsStats "E7" "20191210"
Set-TextValue $wsStats "E8" "67"

# --- Calibration0@0x8 ---
$wsCal = $wb.Worksheets.Item("Calibration0@0x8")
Set-TextValue $wsCal "E3" "-50.440223693847656"
Set-TextValue $wsCal "E9" "-0.03243118152022362"
